# UPDATE hien quan co qua LAN
# Row 13 ("Hiển thị quân cờ đối phương qua mạng LAN") gets its actual
# start/end dates filled in, and the active selection moves to C15.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the actual start/end dates for row 13 (H = actual start, I = actual end)
$ws.Range("H13").Value = 43749
$ws.Range("I13").Value = 43757

# Match the dd/mm date format already used by the other "actual date" cells
$ws.Range("H13:I13").NumberFormat = "dd/mm"

# Move the selection to C15, matching where the user ended up
$ws.Range("C15").Select() | Out-Null
